$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 357, pushing existing rows 357..385 down to 358..386.
$ws.Rows("357:357").Insert()

# Populate the newly inserted row 357 with the new price-record data.
$ws.Range("A357").Value = 4
$ws.Range("B357").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C357").Value = "Los Lagos"
$ws.Range("D357").Value = 44585
$ws.Range("E357").Value = 10
$ws.Range("F357").Value = "Fruta"
$ws.Range("G357").Value = 100102
$ws.Range("H357").Value = "Cítricos"
$ws.Range("I357").Value = 100102003
$ws.Range("J357").Value = "Limón"
$ws.Range("K357").Value = "Sin especificar"
$ws.Range("L357").Value = "1a plateado"
$ws.Range("M357").Value = 800
$ws.Range("N357").Value = 21000
$ws.Range("O357").Value = 21000
$ws.Range("P357").Value = 21000
$ws.Range("Q357").Value = "$/malla 18 kilos"
$ws.Range("R357").Value = "Región de O'Higgins"
$ws.Range("S357").Value = 1167
$ws.Range("T357").Value = 18
